$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.249.61"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "3.144.27"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.144.27"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").Value = "  -2.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.663.97"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.121"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").Value = "63.968.41"
$ws.Range("D18").Value = "3.139.08"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.42"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.99"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "80.92"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +10.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.69"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.11"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "0.0₃0856"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.32"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.06"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.55"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "446.14"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.94"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.288"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0372"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "2.908.95"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.19"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +16.22%  "
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.111"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.17%  "
